$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "65.032.06"
$ws.Cells.Item(2, 4).Style = "Normal"

$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "  +0.95%  "
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.370.16"
$ws.Cells.Item(3, 4).Style = "Normal"

$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "  +0.88%  "
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 4).Style = "Normal"

$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "555.77"
$ws.Cells.Item(5, 4).Style = "Normal"

$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "  +0.36%  "
$ws.Cells.Item(5, 5).Style = "Normal"

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "175.15"
$ws.Cells.Item(6, 4).Style = "Normal"

$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "  -0.06%  "
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "  +2.14%  "
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "3.362.69"
$ws.Cells.Item(8, 4).Style = "Normal"

$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "  +0.87%  "
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "1.00"
$ws.Cells.Item(9, 4).Style = "Normal"

$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.173"
$ws.Cells.Item(10, 4).Style = "Normal"

$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "  +6.31%  "
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.637"
$ws.Cells.Item(11, 4).Style = "Normal"

$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "  +1.74%  "
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "53.78"
$ws.Cells.Item(12, 4).Style = "Normal"

$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "  -1.12%  "
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "  +2.99%  "
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "  +1.23%  "
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.907.20"
$ws.Cells.Item(15, 4).Style = "Normal"

$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "  +0.70%  "
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.120"
$ws.Cells.Item(16, 4).Style = "Normal"

$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "  -0.23%  "
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "3.367.07"
$ws.Cells.Item(18, 4).Style = "Normal"

$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "  +0.69%  "
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "65.040.46"
$ws.Cells.Item(19, 4).Style = "Normal"

$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "  +1.06%  "
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.87"
$ws.Cells.Item(20, 4).Style = "Normal"

$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "  +1.04%  "
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.997"
$ws.Cells.Item(21, 4).Style = "Normal"

$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "  +1.98%  "
$ws.Cells.Item(21, 5).Style = "Normal"

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "452.73"
$ws.Cells.Item(22, 4).Style = "Normal"

$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "  +4.57%  "
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "4.93"
$ws.Cells.Item(23, 4).Style = "Normal"

$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "  -4.22%  "
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "4.08"
$ws.Cells.Item(24, 4).Style = "Normal"

$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "  +0.13%  "
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "86.97"
$ws.Cells.Item(25, 4).Style = "Normal"

$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "  +3.25%  "
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "13.77"
$ws.Cells.Item(26, 4).Style = "Normal"

$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "  +3.02%  "
$ws.Cells.Item(26, 5).Style = "Normal"

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "10.78"
$ws.Cells.Item(27, 4).Style = "Normal"

$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "  +0.32%  "
$ws.Cells.Item(27, 5).Style = "Normal"

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.88"
$ws.Cells.Item(28, 4).Style = "Normal"

$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "  +1.87%  "
$ws.Cells.Item(28, 5).Style = "Normal"

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "8.67"
$ws.Cells.Item(29, 4).Style = "Normal"

$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = "  -0.68%  "
$ws.Cells.Item(29, 5).Style = "Normal"

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "31.06"
$ws.Cells.Item(30, 4).Style = "Normal"

$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = "  +4.60%  "
$ws.Cells.Item(30, 5).Style = "Normal"

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "6.57"
$ws.Cells.Item(31, 4).Style = "Normal"

$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = "  -0.97%  "
$ws.Cells.Item(31, 5).Style = "Normal"

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "62.96"
$ws.Cells.Item(32, 4).Style = "Normal"

$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = "  +8.26%  "
$ws.Cells.Item(32, 5).Style = "Normal"

$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = "  -0.16%  "
$ws.Cells.Item(33, 5).Style = "Normal"

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "576.21"
$ws.Cells.Item(34, 4).Style = "Normal"

$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = "  -0.88%  "
$ws.Cells.Item(34, 5).Style = "Normal"

$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = "  +0.39%  "
$ws.Cells.Item(35, 5).Style = "Normal"

$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = "  +0.08%  "
$ws.Cells.Item(36, 5).Style = "Normal"

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.63"
$ws.Cells.Item(37, 4).Style = "Normal"

$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = "  +4.27%  "
$ws.Cells.Item(37, 5).Style = "Normal"

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.141"
$ws.Cells.Item(38, 4).Style = "Normal"

$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = "  -0.33%  "
$ws.Cells.Item(38, 5).Style = "Normal"

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "35.73"
$ws.Cells.Item(39, 4).Style = "Normal"

$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "  +0.26%  "
$ws.Cells.Item(39, 5).Style = "Normal"

$ws.Cells.Item(40, 2).Value = "TheGraph"

$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.372"
$ws.Cells.Item(40, 4).Style = "Normal"

$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "  +1.31%  "
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 2).Value = "PEPE"

$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0₃0742"
$ws.Cells.Item(41, 4).Style = "Normal"

$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "  -1.05%  "
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.092.26"
$ws.Cells.Item(42, 4).Style = "Normal"

$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "  -0.56%  "
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "  +2.44%  "
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "  -0.75%  "
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "3.19"
$ws.Cells.Item(45, 4).Style = "Normal"

$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "  -0.91%  "
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "  +3.28%  "
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.45"
$ws.Cells.Item(47, 4).Style = "Normal"

$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "  -0.20%  "
$ws.Cells.Item(47, 5).Style = "Normal"

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.00"
$ws.Cells.Item(48, 4).Style = "Normal"

$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "  +0.13%  "
$ws.Cells.Item(48, 5).Style = "Normal"

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "141.59"
$ws.Cells.Item(49, 4).Style = "Normal"

$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = "  +4.91%  "
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.53"
$ws.Cells.Item(50, 4).Style = "Normal"

$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "  -2.30%  "
$ws.Cells.Item(50, 5).Style = "Normal"

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "8.29"
$ws.Cells.Item(51, 4).Style = "Normal"

$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "  +0.43%  "
$ws.Cells.Item(51, 5).Style = "Normal"
